$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.452.60'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.703.13'
$ws.Range('E3').Value = '  +1.06%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.37%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.78'
$ws.Range('E5').Value = '  +0.70%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5479'
$ws.Range('E6').Value = '  +4.32%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.011'
$ws.Range('E7').Value = '  +0.32%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2747'
$ws.Range('E8').Value = '  +1.17%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06472'
$ws.Range('E9').Value = '  +0.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.10'
$ws.Range('E10').Value = '  +0.44%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07712'
$ws.Range('E11').Value = '  +2.93%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.703.67'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.562'
$ws.Range('E13').Value = '  +0.11%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5848'
$ws.Range('E14').Value = '  +0.82%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008406'
$ws.Range('E15').Value = '  -0.69%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.06'
$ws.Range('E16').Value = '  +2.78%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.557.88'
$ws.Range('E17').Value = '  +0.86%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.962'
$ws.Range('E18').Value = '  +0.73%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.011'
$ws.Range('E19').Value = '  +0.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.00'
$ws.Range('E20').Value = '  +1.19%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.45'
$ws.Range('E21').Value = '  +1.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.274'
$ws.Range('E22').Value = '  +1.29%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.012'
$ws.Range('E23').Value = '  +0.26%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.48'
$ws.Range('E24').Value = '  +3.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1327'
$ws.Range('E25').Value = '  +7.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.933'
$ws.Range('E26').Value = '  +2.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.84'
$ws.Range('E27').Value = '  +0.34%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06277'
$ws.Range('E28').Value = '  -5.42%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.380'
$ws.Range('E29').Value = '  +2.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.335'
$ws.Range('E30').Value = '  +0.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.626'
$ws.Range('E31').Value = '  +1.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.609'
$ws.Range('E32').Value = '  +0.94%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.697'
$ws.Range('E33').Value = '  +2.14%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.047'
$ws.Range('E34').Value = '  +1.85%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6200'
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.413'
$ws.Range('E36').Value = '  +0.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('E37').Value = '  +2.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01649'
$ws.Range('E38').Value = '  +1.70%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.122.68'
$ws.Range('E39').Value = '  +1.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.149'
$ws.Range('E40').Value = '  -3.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8845'
$ws.Range('E41').Value = '  +1.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.56'
$ws.Range('E43').Value = '  +0.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.855.78'
$ws.Range('E44').Value = '  +1.23%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('E45').Value = '  +0.15%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.83'
$ws.Range('E46').Value = '  +1.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.247'
$ws.Range('E47').Value = '  +0.71%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  -0.26%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05290'
$ws.Range('E49').Value = '  +0.35%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.167'
$ws.Range('E50').Value = '  +2.18%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4309'
$ws.Range('E51').Value = '  +0.13%  '
